{"js": "// Move the \"LOM3246 - T\u00e9cnicas de Caracteriza\u00e7\u00e3o de Materiais (Indica\u00e7\u00e3o de\n// Conjunto)\" line from the end of the \"Requisitos\" bullet paragraph to the\n// very beginning of that same paragraph (before \"LOB1021 ...\").\n\nconst MOVED_LINE =\n  \"LOM3246 -  T\u00e9cnicas de Caracteriza\u00e7\u00e3o de Materiais  (Indica\u00e7\u00e3o de Conjunto)\";\n\n// Locate the \"Requisitos\" heading paragraph, then the bullet-list paragraph\n// right after it (the one holding the three \"w:br\"-separated requirement\n// lines) \u2014 more robust than a hard-coded paragraph index.\nconst body = context.document.body;\nconst heading = body.search(\"Requisitos\", { matchCase: true, matchWholeWord: true });\nheading.load(\"items\");\nawait context.sync();\n\nif (heading.items.length === 0) {\n  throw new Error(\"Could not find the 'Requisitos' heading paragraph.\");\n}\n\nconst headingParagraph = heading.items[heading.items.length - 1].paragraphs.getFirst();\nconst listParagraph = headingParagraph.getNext();\nlistParagraph.load(\"text\");\nawait context.sync();\n\n// Insert the moved line (with its own trailing manual line break) at the very\n// start of the paragraph. Because the inserted text itself contains the\n// line-break character, Word places it in its own run \u2014 exactly like the\n// other lines already in this paragraph.\nconst target = listParagraph.getRange();\ntarget.insertText(MOVED_LINE + \"\\u000b\", \"Start\");\nawait context.sync();\n\n// The paragraph now contains the moved line twice: the fresh copy at the\n// start, and the original one still sitting at the end. Find every\n// occurrence (including its trailing break so the match spans the whole\n// former run) and drop the last one \u2014 the original \u2014 leaving only the new\n// copy up front.\nconst dupes = listParagraph.search(MOVED_LINE + \"\\u000b\", { matchCase: true });\ndupes.load(\"items\");\nawait context.sync();\n\nif (dupes.items.length < 2) {\n  throw new Error(\"Expected to find the duplicated requirement line twice.\");\n}\n\ndupes.items[dupes.items.length - 1].delete();\nawait context.sync();\n", "ps1": "# Move the \"LOM3246 - T\u00e9cnicas de Caracteriza\u00e7\u00e3o de Materiais (Indica\u00e7\u00e3o de\n# Conjunto)\" line from the end of the \"Requisitos\" bullet paragraph to the\n# very beginning of that same paragraph (before \"LOB1021 ...\").\n\n$d = $word.ActiveDocument\n$movedLine = \"LOM3246 -  T\u00e9cnicas de Caracteriza\u00e7\u00e3o de Materiais  (Indica\u00e7\u00e3o de Conjunto)\"\n\n# Locate the \"Requisitos\" heading paragraph via Find (robust to a hard-coded\n# paragraph index), then grab the very next paragraph \u2014 the bullet list that\n# holds the three w:br-separated requirement lines.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Requisitos\")\nif (-not $found) {\n    throw \"Could not find the 'Requisitos' heading paragraph.\"\n}\n\n$precedingRange = $d.Range(0, $findRange.Start)\n$headingIndex = $precedingRange.Paragraphs.Count + 1\n$listParagraph = $d.Paragraphs.Item($headingIndex + 1)\n$listRange = $listParagraph.Range\n\n# Insert the moved line (with its own trailing manual line break) at the very\n# start of the paragraph. Because the inserted text itself contains the\n# line-break character, Word places it in its own run - exactly like the\n# other lines already in this paragraph.\n$listRange.InsertBefore($movedLine + [char]11)\n\n# The paragraph now contains the moved line twice: the fresh copy at the\n# start, and the original one still sitting at the end. Search the remainder\n# of the paragraph (everything after our freshly inserted copy) for the next\n# occurrence, extend the match by one character to also swallow its trailing\n# line break, and delete it - leaving only the new copy up front.\n$searchStart = $listRange.Start + $movedLine.Length + 1\n$searchRange = $d.Range($searchStart, $listParagraph.Range.End)\n$dupFound = $searchRange.Find.Execute($movedLine)\nif (-not $dupFound) {\n    throw \"Could not find the duplicated requirement line.\"\n}\n$searchRange.MoveEnd(1, 1)\n$searchRange.Delete()\n"}
